$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Вопросы" (sheet1): insert a new column B "Номер_Вопроса" holding
# the question number (1,2,3). The existing "Вопрос" and "Дата создания"
# columns shift right by one (to C and D respectively).
# ----------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Вопросы")
$wsQ.Columns("B").Insert()

# Clear the formatting that Insert() copied from column A so the new
# number cells stay unstyled, matching the target layout.
$wsQ.Range("B2:B4").Clear()

$wsQ.Range("B2").Value = 1
$wsQ.Range("B3").Value = 2
$wsQ.Range("B4").Value = 3

# Copy header formatting (bold / bordered style) from the neighboring
# header cell onto the new header cell, then set its text.
$wsQ.Range("C1").Copy()
$wsQ.Range("B1").PasteSpecial(-4122)
$wsQ.Range("B1").Value = "Номер_Вопроса"

# ----------------------------------------------------------------------
# Sheet "Ответы" (sheet2): update unique_key (B), Оценка (D) and ПВИ (G)
# values for several groups of rows as part of the PVI drop-down rework.
# ----------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Ответы")

$wsA.Range("B2").Value = 179182
$wsA.Range("B3").Value = 179182
$wsA.Range("B4").Value = 179182
$wsA.Range("G2").Value = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)"
$wsA.Range("G3").Value = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)"
$wsA.Range("G4").Value = "Н-Уренгойское ЛПУМГ (ПВП № 1, КС Пуртазовская)"

$wsA.Range("B5").Value = 212300
$wsA.Range("B6").Value = 212300
$wsA.Range("B7").Value = 212300
$wsA.Range("D5").Value = 1
$wsA.Range("D6").Value = 1
$wsA.Range("D7").Value = 1
$wsA.Range("G5").Value = "Н-Уренгойское ЛПУМГ (ПВП №2)"
$wsA.Range("G6").Value = "Н-Уренгойское ЛПУМГ (ПВП №2)"
$wsA.Range("G7").Value = "Н-Уренгойское ЛПУМГ (ПВП №2)"

$wsA.Range("B8").Value = 13679
$wsA.Range("B9").Value = 13679
$wsA.Range("B10").Value = 13679
$wsA.Range("D8").Value = 1
$wsA.Range("D9").Value = 1
$wsA.Range("D10").Value = 1
$wsA.Range("G8").Value = "Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)"
$wsA.Range("G9").Value = "Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)"
$wsA.Range("G10").Value = "Пурпейское ЛПУМГ (КС - 01, Общежитие на 100 мест Ягенетская п/п)"

$wsA.Range("B11").Value = 146533
$wsA.Range("B12").Value = 146533
$wsA.Range("B13").Value = 146533
